$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mapping of row -> (code in column A, expected old hash, new hash) for column B
# Verified against the source workbook; each update is applied only if the
# existing value matches the expected old hash, to guard against drift.
$updates = @(
    @{ Row = 54; Code = "01-010063TM"; Old = "bf44ec96ab80cb1716583fd2713be6b3"; New = "c158a919d379f2b0a1dd3f92c64e0efe" }
    @{ Row = 80; Code = "01-010063TC"; Old = "7cd71806c3817a2788b411cc5dc0d07f"; New = "90b938ed5724b1518f61528df0964e9b" }
    @{ Row = 108; Code = "01-010065TM"; Old = "1c6b965ee60990ab7717581b4a83445d"; New = "9d9599fe4b082dacd3c78b3092d2d689" }
    @{ Row = 159; Code = "05-050203TP"; Old = "5516c0461909e150764e8279d36584a4"; New = "14b3561a331a489dc600983a00a250ea" }
    @{ Row = 216; Code = "03-030077TC"; Old = "f1173a4f1a4db9a6ccec48b3aff58281"; New = "0e404721817726b5a49708e58f8cbca6" }
    @{ Row = 225; Code = "03-030077TP"; Old = "0356828e5f648d949da58190375d702e"; New = "46ca80ee38557fe7b8ae34646f1f9b04" }
    @{ Row = 246; Code = "05-050003TC"; Old = "a7844963b70be534ed450364d9f7d1e9"; New = "939d2c7d7063e6792a0325ce3de3605b" }
    @{ Row = 276; Code = "05-050003TP"; Old = "a5a8399642eb3856bc0ed3d26c605c8e"; New = "650806b026dded975bb0cf915c895d82" }
    @{ Row = 281; Code = "05-050201TC"; Old = "17cfa0728bacabad7c7d2276ad59d422"; New = "bb5b6bee92418c56be72f11c3683224e" }
    @{ Row = 330; Code = "05-050005TC"; Old = "641c214c0ff497c231e16e0202107c57"; New = "1d179338b2bef50a133c6811b4d220d3" }
    @{ Row = 339; Code = "05-050201TP"; Old = "bb925f9bce4146dbc18f0ef0f1387cf2"; New = "11cc37d6882ae63c07d8776d0f0d0ccc" }
    @{ Row = 419; Code = "05-0709-070905BTC"; Old = "afba4ee92bb44bede48ddf483ac24705"; New = "2f36e7fae61a39e97cd825cd8a551d49" }
    @{ Row = 446; Code = "05-050001TP"; Old = "9de5a67740a3686774a6f39010a19265"; New = "c31f4f6d58b96caebcf7986a5ea14961" }
    @{ Row = 460; Code = "05-050204A"; Old = "a5dbe54c39a9069dfff780add106e62d"; New = "18a56c4539a9a8fe0481b0ccd34dd7fc" }
    @{ Row = 500; Code = "05-050202A"; Old = "699658c5c4dee4e8bbbd60f12d5ecc22"; New = "10f73d38869a48511c01da7a4ba5a258" }
    @{ Row = 517; Code = "05-050203A"; Old = "790260fc1b06d5c1e5750256043dad45"; New = "80869109392a13b261410e0fcad571f4" }
    @{ Row = 543; Code = "01-010064TM"; Old = "ece6eb734faed0dd6d9b51a279f5053d"; New = "ecf2e8e7c4e062d14712d31394abc565" }
    @{ Row = 566; Code = "05-050004A"; Old = "93cf8370596863b200b01bd187da9d14"; New = "26d728b2387296374b27bcac23a3eeaa" }
    @{ Row = 574; Code = "01-010064TC"; Old = "a129a870088d76f781fe1f5950d3a8ba"; New = "94b75164276b0a69d6415771a737a4cc" }
    @{ Row = 575; Code = "01-010064A"; Old = "03f38022c575245c28fc04992de3c384"; New = "13d93fb487c697baef92569706486f45" }
    @{ Row = 601; Code = "03-030078TC"; Old = "a9ea093c40eaf3e1f00e4a1907276733"; New = "122e50541cdb47f369c40eb3484e3e6d" }
    @{ Row = 606; Code = "03-030078TM"; Old = "7f7ab1f8dc3ebc7cf76fcb6d6f79cd33"; New = "ddfc15fa35a7ccbb61c1e0a1b8fbc20f" }
    @{ Row = 614; Code = "03-030078TP"; Old = "7196dec3ea8c8be2c644d2ff1202802c"; New = "66b5b0a08edb99e212cc7a3f8ad0be3c" }
    @{ Row = 616; Code = "05-050204TP"; Old = "0a8197a280321a7f99dd9c791f024dce"; New = "081ef3448a92487cebead750118c75db" }
    @{ Row = 627; Code = "05-050204TC"; Old = "2521330e9c43a86a2061c5c26fcd442a"; New = "7f36f63d4a8b866cffb76db4148f8a0c" }
    @{ Row = 756; Code = "05-050004TC"; Old = "9397a483900340432a332a438b43feee"; New = "a833e19224e0d52ed7ff59b2093d743d" }
    @{ Row = 761; Code = "05-050004TP"; Old = "9986aac1f2a947465545084339a92eed"; New = "8943035717e9dd7224e03d1d65866774" }
    @{ Row = 786; Code = "05-050002TC"; Old = "d7c32f6feaa74b68ad82f3fb3036d04e"; New = "cce1431dbcb001b2f9256b5a751c030a" }
    @{ Row = 850; Code = "01-010063A"; Old = "ee5f9b6f034b61262ef8922f4d4f5ebd"; New = "295419e181064dbd29a76f9d8800925e" }
    @{ Row = 855; Code = "05-050002A"; Old = "ec5110340224ff40e879ea2857e85751"; New = "bf7ab094b4c95ef6d20396c4daf40bbd" }
    @{ Row = 869; Code = "05-050003A"; Old = "87d5f4401301379682bc0ad75b7a1ef8"; New = "0fa58d77443080fe9884b2d255891018" }
    @{ Row = 904; Code = "05-050001A"; Old = "cd1a090fd82a983cf3eef5f74f74fdd1"; New = "9f4b2b72827258dde43dd35f9a7136ec" }
    @{ Row = 914; Code = "03-030078A"; Old = "012a3efc3a13ac4e2a1886c163d35e01"; New = "1a08015e8770bd1d5e9d50d26eac86b0" }
    @{ Row = 925; Code = "03-030077A"; Old = "3f27e8aadd43ec8a51d6e3542f7dce0e"; New = "34e6d2e0525c35113cab9db343c6672e" }
    @{ Row = 928; Code = "05-050002TP"; Old = "075dc0b3177c298bc5836ccf2890df11"; New = "14d518e206a68cad904790671a1cdd61" }
)

foreach ($u in $updates) {
    $codeCell = $ws.Cells.Item($u.Row, 1)
    $hashCell = $ws.Cells.Item($u.Row, 2)
    if ($codeCell.Value2 -eq $u.Code -and $hashCell.Value2 -eq $u.Old) {
        $hashCell.Value = $u.New
    } else {
        Write-Output "Mismatch at row $($u.Row): code=$($codeCell.Value2) hash=$($hashCell.Value2)"
    }
}
